# Racquet inventory workbook update:
#  - Bold the header row (A1:I1)
#  - Remove the old (empty/unused) DESCRIPTION column J
#  - Fix a few stat typos in the existing rows (I2, H3/I3, I4)
#  - Add two new racquets (rows 5 and 6)
#  - Set the sheet up for printing (portrait page setup) and update the
#    selected cell to reflect where the user left off

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: bold formatting ---------------------------------------
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1:I1").Font.Bold = $true

# --- Drop the unused DESCRIPTION column (J) -----------------------------
$ws.Range("J1").ClearContents()
$ws.Columns(10).Delete()

# --- Correct a few existing data points ---------------------------------
$ws.Range("I2").Value = 3

$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 4

$ws.Range("I4").Value = 3

# --- Add two new racquets -------------------------------------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Mizuno"
$ws.Range("C5").Value = "Fortius 30 Power"
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 5

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Carlton"
$ws.Range("C6").Value = "Air Edge"
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 2

# --- Page setup for printing ---------------------------------------------
$ws.PageSetup.Orientation = 1

# --- Leave the selection where the user ended up --------------------------
$ws.Range("J6").Select()
